{"js": "// Replace the date line and each of the 25 division problems with their\n// updated values, per the commit's diff. Every source string is unique in\n// the document, so a simple search + full-text replace per pair is safe.\nconst replacements = [\n  [\"2024-12-17 Tuesday\", \"2024-12-18 Wednesday\"],\n  [\"415\u00f77=\", \"339\u00f79=\"],\n  [\"936\u00f79=\", \"766\u00f73=\"],\n  [\"834\u00f74=\", \"355\u00f75=\"],\n  [\"449\u00f75=\", \"713\u00f79=\"],\n  [\"244\u00f72=\", \"399\u00f77=\"],\n  [\"177\u00f73=\", \"927\u00f72=\"],\n  [\"442\u00f77=\", \"354\u00f78=\"],\n  [\"233\u00f74=\", \"301\u00f77=\"],\n  [\"893\u00f73=\", \"105\u00f73=\"],\n  [\"601\u00f75=\", \"813\u00f77=\"],\n  [\"341\u00f73=\", \"157\u00f77=\"],\n  [\"291\u00f72=\", \"606\u00f72=\"],\n  [\"382\u00f77=\", \"134\u00f77=\"],\n  [\"602\u00f73=\", \"141\u00f76=\"],\n  [\"605\u00f73=\", \"453\u00f77=\"],\n  [\"218\u00f77=\", \"623\u00f76=\"],\n  [\"452\u00f77=\", \"759\u00f79=\"],\n  [\"107\u00f77=\", \"661\u00f76=\"],\n  [\"639\u00f75=\", \"887\u00f79=\"],\n  [\"456\u00f73=\", \"107\u00f74=\"],\n  [\"504\u00f74=\", \"289\u00f77=\"],\n  [\"292\u00f78=\", \"881\u00f73=\"],\n  [\"406\u00f74=\", \"295\u00f78=\"],\n  [\"913\u00f76=\", \"620\u00f79=\"],\n  [\"420\u00f76=\", \"113\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each of the 25 division problems with their\n# updated values, per the commit's diff. Every source string is unique in\n# the document, so Find/Replace (wdReplaceAll) per pair is safe and exact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-12-17 Tuesday\", \"2024-12-18 Wednesday\"),\n  @(\"415\u00f77=\", \"339\u00f79=\"),\n  @(\"936\u00f79=\", \"766\u00f73=\"),\n  @(\"834\u00f74=\", \"355\u00f75=\"),\n  @(\"449\u00f75=\", \"713\u00f79=\"),\n  @(\"244\u00f72=\", \"399\u00f77=\"),\n  @(\"177\u00f73=\", \"927\u00f72=\"),\n  @(\"442\u00f77=\", \"354\u00f78=\"),\n  @(\"233\u00f74=\", \"301\u00f77=\"),\n  @(\"893\u00f73=\", \"105\u00f73=\"),\n  @(\"601\u00f75=\", \"813\u00f77=\"),\n  @(\"341\u00f73=\", \"157\u00f77=\"),\n  @(\"291\u00f72=\", \"606\u00f72=\"),\n  @(\"382\u00f77=\", \"134\u00f77=\"),\n  @(\"602\u00f73=\", \"141\u00f76=\"),\n  @(\"605\u00f73=\", \"453\u00f77=\"),\n  @(\"218\u00f77=\", \"623\u00f76=\"),\n  @(\"452\u00f77=\", \"759\u00f79=\"),\n  @(\"107\u00f77=\", \"661\u00f76=\"),\n  @(\"639\u00f75=\", \"887\u00f79=\"),\n  @(\"456\u00f73=\", \"107\u00f74=\"),\n  @(\"504\u00f74=\", \"289\u00f77=\"),\n  @(\"292\u00f78=\", \"881\u00f73=\"),\n  @(\"406\u00f74=\", \"295\u00f78=\"),\n  @(\"913\u00f76=\", \"620\u00f79=\"),\n  @(\"420\u00f76=\", \"113\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
